# Microcontroller Peripheral Map - ADC driver update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (UART3): software naming renamed USB_UART -> usbUart ---
$ws.Range("D5").Value = "usbUart"

# --- Row 7 (TMR1): note no longer mentions clearing the watchdog timer ---
$ws.Range("C7").Value = "Timer 1 is used as the heartbeat timer, which triggers an interrupt every second. This blinks an LED and increments on time counters. Used as a 1 Hz timebase"

# --- New row 8: TMR3 / Timer 3 used for ADC trigger timebase ---
$ws.Range("A8").Value = "TMR3"
$ws.Range("B8").Value = "Timer 3"
$ws.Range("C8").Value = "Used for timebase for automatic ADC triggering"
$ws.Range("D8").Value = "ADCTriggerTimer"

# --- New row 9: DMT / Deadman timer ---
$ws.Range("A9").Value = "DMT"
$ws.Range("B9").Value = "Deadman timer"
$ws.Range("C9").Value = "Resets the microcontroller if not cleared by the core timer ISR"
$ws.Range("D9").Value = "Deadman Timer"

# --- New row 10: Core Timer ---
$ws.Range("A10").Value = "Core Timer"
$ws.Range("B10").Value = "Core Timer"
$ws.Range("C10").Value = "CP0 integrated core timer, used for clearing the WDT and DMT"
$ws.Range("D10").Value = "Coretimer"

# --- Row heights ---
$ws.Rows.Item(7).RowHeight = 72.5
$ws.Rows.Item(8).RowHeight = 29
$ws.Rows.Item(9).RowHeight = 29
$ws.Rows.Item(10).RowHeight = 29

# --- Selection moves to D2 ---
$ws.Range("D2").Select() | Out-Null
